$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at E:F (pushes old E:J -> G:L)
$ws.Range("E1:F1").EntireColumn.Insert()

# New header cells
$ws.Range("E1").Value = "Phone"
$ws.Range("F1").Value = "WhatsApp Enabled"

# New data cells
$ws.Range("E2").Value = 9999999999
$ws.Range("F2").Value = "Yes"

$ws.Range("E3").Value = 8888888888
$ws.Range("F3").Value = "No"

# Update selection to match target
$ws.Range("F4").Select() | Out-Null
